# Rename the three municipal office destinations in the "rutas" table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "concejo"
$ws.Range("B5").Value = "despacho"
$ws.Range("B6").Value = "cobro coactivo"

# Update the active selection to C7.
$ws.Range("C7").Select()
